$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly-measured ATtinyX4 (E/F) and Total (G/H) columns ---
# Row 3
$ws.Range("E3").Value = 1446
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1684
$ws.Range("H3").Value = 1

# Row 4
$ws.Range("E4").Value = 954
$ws.Range("F4").Value = 104
$ws.Range("G4").Value = 1132
$ws.Range("H4").Value = 104

# Row 5
$ws.Range("E5").Value = 1232
$ws.Range("F5").Value = 104
$ws.Range("G5").Value = 1410
$ws.Range("H5").Value = 104

# Row 6
$ws.Range("E6").Value = 1238
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 1416
$ws.Range("H6").Value = 104

# Row 7
$ws.Range("E7").Value = 310
$ws.Range("F7").Value = 7
$ws.Range("G7").Value = 522
$ws.Range("H7").Value = 7

# Row 8 - code size for ATmega328P changed, and the old value moved to the
# (now-unused for this row) old ATmega2560 columns J/K
$ws.Range("C8").Value = 156
$ws.Range("E8").Value = 86
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 284
$ws.Range("H8").Value = 0
$ws.Range("C8").Copy() | Out-Null
$ws.Range("J8:K8").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$ws.Range("J8").Value = 174
$ws.Range("K8").Value = 0

# Row 10
$ws.Range("E10").Value = 162
$ws.Range("F10").Value = 0

# Row 11
$ws.Range("E11").Value = 276
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 474
$ws.Range("H11").Value = 0

# Row 12 - no ATtinyX4 build for this example ("-" placeholder, right-aligned)
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "-"
$ws.Range("E12:F12").HorizontalAlignment = -4152 # xlRight
$ws.Range("G12").Value = 2031
$ws.Range("H12").Value = 169

# Row 13
$ws.Range("E13").Value = 1584
$ws.Range("F13").Value = 158
$ws.Range("G13").Value = 1768
$ws.Range("H13").Value = 158

# Row 14
$ws.Range("E14").Value = 1582
$ws.Range("F14").Value = 158
$ws.Range("G14").Value = 1766
$ws.Range("H14").Value = 158

# Row 15 - new description + measurements
$ws.Range("B15").Value = "Use PCI on push button to light LED"
$ws.Range("E15").Value = 430
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = 652
$ws.Range("H15").Value = 8

# Row 16 - new description + measurements
$ws.Range("B16").Value = "Use PCI on 3 push buttons to light 4 LEDs"
$ws.Range("E16").Value = 618
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 820
$ws.Range("H16").Value = 8

# Row 17 - new description + measurements
$ws.Range("B17").Value = "Use PCI on 3 push buttons to light 4 LEDs"
$ws.Range("E17").Value = 408
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 626
$ws.Range("H17").Value = 8

# --- Extend the (now unused for most rows) ATmega2560 J/K columns so every
# data row keeps matching formatting all the way across, even where there is
# no longer a distinct ATmega2560 measurement ---
$ws.Range("C7").Copy() | Out-Null
$ws.Range("J7:K7").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$ws.Range("C9").Copy() | Out-Null
$ws.Range("J9:K9").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("J10:K10").PasteSpecial(-4122) | Out-Null
$ws.Range("C11").Copy() | Out-Null
$ws.Range("J11:K11").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("J13:K13").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("J14:K14").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Copy() | Out-Null
$ws.Range("J15:K15").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Copy() | Out-Null
$ws.Range("J16:K16").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").Copy() | Out-Null
$ws.Range("J17:K17").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Restore selection to where the author left off ---
$ws.Range("F17").Select() | Out-Null
